# Scheduled market-price refresh: update currentAveragePrice(/NQ/HQ),
# LevePriceNQ/HQ and LeveProfitNQ/HQ columns (H,I,J,K,L,M,N) for the rows
# whose Universalis price snapshot changed, across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 497.75
$ws.Range("I9").Value = 430.33334
$ws.Range("K9").Value = 430.33334
$ws.Range("M9").Value = -261.33334

$ws.Range("H19").Value = 4994.5
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H29").Value = 3938.125
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 4357.857
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 13073.571
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = -13635.571

$ws.Range("H33").Value = 2050.423
$ws.Range("I33").Value = 1628.4
$ws.Range("J33").Value = 3457.1667
$ws.Range("K33").Value = 1628.4
$ws.Range("L33").Value = 3457.1667
$ws.Range("M33").Value = -1399.4
$ws.Range("N33").Value = -3915.1667

$ws.Range("H52").Value = 1377.9
$ws.Range("I52").Value = 420
$ws.Range("K52").Value = 1260
$ws.Range("M52").Value = -1100

$ws.Range("H63").Value = 113000
$ws.Range("J63").Value = 113000
$ws.Range("L63").Value = 113000
$ws.Range("N63").Value = -114248

$ws.Range("H66").Value = 113000
$ws.Range("J66").Value = 113000
$ws.Range("L66").Value = 339000
$ws.Range("N66").Value = -345240

$ws.Range("H135").Value = 2897.75
$ws.Range("I135").Value = 2757.6
$ws.Range("K135").Value = 24818.4
$ws.Range("M135").Value = -22283.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17865602
$ws.Range("I32").Value = 20836868
$ws.Range("K32").Value = 20836868
$ws.Range("M32").Value = -20836581

$ws.Range("H61").Value = 37507348
$ws.Range("I61").Value = 27784708
$ws.Range("K61").Value = 27784708
$ws.Range("M61").Value = -27784496

$ws.Range("H63").Value = 5464.7144
$ws.Range("I63").Value = 2357.1428
$ws.Range("K63").Value = 2357.1428
$ws.Range("M63").Value = -1671.1428

$ws.Range("H66").Value = 5464.7144
$ws.Range("I66").Value = 2357.1428
$ws.Range("K66").Value = 11785.714
$ws.Range("M66").Value = -8353.714

$ws.Range("H136").Value = 37507348
$ws.Range("I136").Value = 27784708
$ws.Range("K136").Value = 83354124
$ws.Range("M136").Value = -83351574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 4444
$ws.Range("J45").Value = 4444
$ws.Range("L45").Value = 4444
$ws.Range("N45").Value = -6060

$ws.Range("H76").Value = 53313.332
$ws.Range("J76").Value = 57978
$ws.Range("L76").Value = 57978
$ws.Range("N76").Value = -58608

$ws.Range("H79").Value = 53313.332
$ws.Range("J79").Value = 57978
$ws.Range("L79").Value = 57978
$ws.Range("N79").Value = -60162

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 451.8889
$ws.Range("I22").Value = 380
$ws.Range("K22").Value = 380
$ws.Range("M22").Value = -30

$ws.Range("H58").Value = 2674.3044
$ws.Range("I58").Value = 1458.421
$ws.Range("K58").Value = 1458.421
$ws.Range("M58").Value = -1255.421

$ws.Range("H107").Value = 2078.7
$ws.Range("I107").Value = 895.8570999999999
$ws.Range("K107").Value = 895.8570999999999
$ws.Range("M107").Value = 1024.1429

$ws.Range("H132").Value = 1776.2307
$ws.Range("I132").Value = 1799.25
$ws.Range("K132").Value = 5397.75
$ws.Range("M132").Value = -2867.75

$ws.Range("H136").Value = 2674.3044
$ws.Range("I136").Value = 1458.421
$ws.Range("K136").Value = 4375.263
$ws.Range("M136").Value = -1825.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 756.5714
$ws.Range("I6").Value = 716
$ws.Range("K6").Value = 2148
$ws.Range("M6").Value = -2035

$ws.Range("H25").Value = 211.11765
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H29").Value = 92
$ws.Range("I29").Value = 53
$ws.Range("J29").Value = 119.85714
$ws.Range("K29").Value = 159
$ws.Range("L29").Value = 359.57142
$ws.Range("M29").Value = 118
$ws.Range("N29").Value = -913.57142

$ws.Range("H30").Value = 211.11765
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H35").Value = 299.5
$ws.Range("J35").Value = 299.5
$ws.Range("L35").Value = 898.5
$ws.Range("N35").Value = -1474.5

$ws.Range("H44").Value = 558.6
$ws.Range("I44").Value = 323.25
$ws.Range("J44").Value = 1500
$ws.Range("K44").Value = 969.75
$ws.Range("L44").Value = 4500
$ws.Range("M44").Value = -571.75
$ws.Range("N44").Value = -5296

$ws.Range("H47").Value = 3741.3
$ws.Range("J47").Value = 12499.5
$ws.Range("L47").Value = 37498.5
$ws.Range("N47").Value = -38360.5

$ws.Range("H113").Value = 1228.3077
$ws.Range("J113").Value = 1315.3636
$ws.Range("L113").Value = 3946.0908
$ws.Range("N113").Value = -8286.0908

$ws.Range("H121").Value = 842.46155
$ws.Range("J121").Value = 912.25
$ws.Range("L121").Value = 2736.75
$ws.Range("N121").Value = -5356.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6280.727
$ws.Range("I70").Value = 5038
$ws.Range("J70").Value = 7316.3335
$ws.Range("K70").Value = 5038
$ws.Range("L70").Value = 7316.3335
$ws.Range("M70").Value = -4768
$ws.Range("N70").Value = -7856.3335

$ws.Range("H73").Value = 6280.727
$ws.Range("I73").Value = 5038
$ws.Range("J73").Value = 7316.3335
$ws.Range("K73").Value = 5038
$ws.Range("L73").Value = 7316.3335
$ws.Range("M73").Value = -4102
$ws.Range("N73").Value = -9188.333500000001

$ws.Range("H132").Value = 24396766
$ws.Range("J132").Value = 15249.444
$ws.Range("L132").Value = 45748.33199999999
$ws.Range("N132").Value = -50808.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3137.2856
$ws.Range("I100").Value = 2115.25
$ws.Range("K100").Value = 2115.25
$ws.Range("M100").Value = -1574.25

$ws.Range("H136").Value = 101675.08
$ws.Range("I136").Value = 14271.375
$ws.Range("K136").Value = 42814.125
$ws.Range("M136").Value = -40264.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 21666.666
$ws.Range("J69").Value = 21666.666
$ws.Range("L69").Value = 21666.666
$ws.Range("N69").Value = -23164.666

$ws.Range("H72").Value = 21666.666
$ws.Range("J72").Value = 21666.666
$ws.Range("L72").Value = 64999.99800000001
$ws.Range("N72").Value = -72487.99800000001

$ws.Range("H136").Value = 3845
$ws.Range("I136").Value = 3737.25
$ws.Range("K136").Value = 11211.75
$ws.Range("M136").Value = -8661.75
